$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2275.2563
$ws.Range("I51").Value = 1800.6818
$ws.Range("J51").Value = 2889.4119
$ws.Range("K51").Value = 1800.6818
$ws.Range("L51").Value = 2889.4119
$ws.Range("M51").Value = -1316.6818
$ws.Range("N51").Value = -3857.4119
$ws.Range("H74").Value = 5199.3335
$ws.Range("I74").Value = 3598
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 3598
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -2662
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 5199.3335
$ws.Range("I77").Value = 3598
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 17990
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -13310
$ws.Range("N77").Value = -39360
$ws.Range("H100").Value = 1672.6
$ws.Range("I100").Value = 768.0714
$ws.Range("J100").Value = 3783.1667
$ws.Range("K100").Value = 768.0714
$ws.Range("L100").Value = 3783.1667
$ws.Range("M100").Value = -227.0714
$ws.Range("N100").Value = -4865.1667
$ws.Range("H103").Value = 684.6667
$ws.Range("J103").Value = 684.6667
$ws.Range("L103").Value = 2054.0001
$ws.Range("N103").Value = -3226.0001
$ws.Range("H137").Value = 66669224
$ws.Range("I137").Value = 38464180
$ws.Range("J137").Value = 250002000
$ws.Range("K137").Value = 115392540
$ws.Range("L137").Value = 750006000
$ws.Range("M137").Value = -115389990
$ws.Range("N137").Value = -750011100
$ws.Range("H138").Value = 3291.17
$ws.Range("I138").Value = 2036.3889
$ws.Range("J138").Value = 3936.4856
$ws.Range("K138").Value = 6109.1667
$ws.Range("L138").Value = 11809.4568
$ws.Range("M138").Value = -969.1666999999998
$ws.Range("N138").Value = -22089.4568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1853.875
$ws.Range("I45").Value = 1670.75
$ws.Range("J45").Value = 2037
$ws.Range("K45").Value = 1670.75
$ws.Range("L45").Value = 2037
$ws.Range("M45").Value = -1293.75
$ws.Range("N45").Value = -2791
$ws.Range("H61").Value = 35718200
$ws.Range("I61").Value = 47621796
$ws.Range("K61").Value = 47621796
$ws.Range("M61").Value = -47621584
$ws.Range("H132").Value = 26324368
$ws.Range("I132").Value = 8719.5
$ws.Range("J132").Value = 166674500
$ws.Range("K132").Value = 26158.5
$ws.Range("L132").Value = 500023500
$ws.Range("M132").Value = -23628.5
$ws.Range("N132").Value = -500028560
$ws.Range("H136").Value = 35718200
$ws.Range("I136").Value = 47621796
$ws.Range("K136").Value = 142865388
$ws.Range("M136").Value = -142862838

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 20837.055
$ws.Range("I86").Value = 12620.546
$ws.Range("J86").Value = 33748.715
$ws.Range("K86").Value = 12620.546
$ws.Range("L86").Value = 33748.715
$ws.Range("M86").Value = -11497.546
$ws.Range("N86").Value = -35994.715
$ws.Range("H89").Value = 20837.055
$ws.Range("I89").Value = 12620.546
$ws.Range("J89").Value = 33748.715
$ws.Range("K89").Value = 63102.73
$ws.Range("L89").Value = 168743.575
$ws.Range("M89").Value = -57486.73
$ws.Range("N89").Value = -179975.575
$ws.Range("H134").Value = 4856.8213
$ws.Range("I134").Value = 3045.9546
$ws.Range("K134").Value = 9137.863799999999
$ws.Range("M134").Value = -6602.863799999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33338474
$ws.Range("I31").Value = 4229.409
$ws.Range("K31").Value = 4229.409
$ws.Range("M31").Value = -3934.409
$ws.Range("H34").Value = 33338474
$ws.Range("I34").Value = 4229.409
$ws.Range("K34").Value = 4229.409
$ws.Range("M34").Value = -4027.409
$ws.Range("H58").Value = 4861.357
$ws.Range("I58").Value = 4896.5835
$ws.Range("J58").Value = 4650
$ws.Range("K58").Value = 4896.5835
$ws.Range("L58").Value = 4650
$ws.Range("M58").Value = -4693.5835
$ws.Range("N58").Value = -5056
$ws.Range("H99").Value = 6029.96
$ws.Range("I99").Value = 6687.1055
$ws.Range("K99").Value = 6687.1055
$ws.Range("M99").Value = -5189.1055
$ws.Range("H126").Value = 6029.96
$ws.Range("I126").Value = 6687.1055
$ws.Range("K126").Value = 20061.3165
$ws.Range("M126").Value = -17591.3165
$ws.Range("H132").Value = 85898.25
$ws.Range("I132").Value = 101916.5
$ws.Range("J132").Value = 5807
$ws.Range("K132").Value = 305749.5
$ws.Range("L132").Value = 17421
$ws.Range("M132").Value = -303219.5
$ws.Range("N132").Value = -22481
$ws.Range("H136").Value = 4861.357
$ws.Range("I136").Value = 4896.5835
$ws.Range("J136").Value = 4650
$ws.Range("K136").Value = 14689.7505
$ws.Range("L136").Value = 13950
$ws.Range("M136").Value = -12139.7505
$ws.Range("N136").Value = -19050

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2089.8667
$ws.Range("I109").Value = 1394.3334
$ws.Range("K109").Value = 4183.0002
$ws.Range("M109").Value = -3143.0002
$ws.Range("H117").Value = 3324.818
$ws.Range("I117").Value = 1606.3334
$ws.Range("J117").Value = 3969.25
$ws.Range("K117").Value = 4819.0002
$ws.Range("L117").Value = 11907.75
$ws.Range("M117").Value = -1377.0002
$ws.Range("N117").Value = -18791.75
$ws.Range("H121").Value = 3343.3333
$ws.Range("I121").Value = 4030
$ws.Range("K121").Value = 12090
$ws.Range("M121").Value = -10780

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4427.769
$ws.Range("J80").Value = 4446.1
$ws.Range("L80").Value = 4446.1
$ws.Range("N80").Value = -6442.1
$ws.Range("H83").Value = 4427.769
$ws.Range("J83").Value = 4446.1
$ws.Range("L83").Value = 22230.5
$ws.Range("N83").Value = -32214.5
$ws.Range("H102").Value = 1629.125
$ws.Range("I102").Value = 1394.35
$ws.Range("J102").Value = 2803
$ws.Range("K102").Value = 1394.35
$ws.Range("L102").Value = 2803
$ws.Range("M102").Value = 227.6500000000001
$ws.Range("N102").Value = -6047
$ws.Range("H132").Value = 2068.196
$ws.Range("I132").Value = 1859.7826
$ws.Range("K132").Value = 5579.3478
$ws.Range("M132").Value = -3049.3478

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4098.3887
$ws.Range("I7").Value = 3720.8076
$ws.Range("K7").Value = 3720.8076
$ws.Range("M7").Value = -3608.8076
$ws.Range("H114").Value = 52198
$ws.Range("J114").Value = 52198
$ws.Range("L114").Value = 52198
$ws.Range("N114").Value = -60876
$ws.Range("H126").Value = 4098.3887
$ws.Range("I126").Value = 3720.8076
$ws.Range("K126").Value = 11162.4228
$ws.Range("M126").Value = -8692.4228
$ws.Range("H136").Value = 3174.9756
$ws.Range("I136").Value = 3231.875
$ws.Range("J136").Value = 899
$ws.Range("K136").Value = 9695.625
$ws.Range("L136").Value = 2697
$ws.Range("M136").Value = -7145.625
$ws.Range("N136").Value = -7797

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 44500
$ws.Range("J116").Value = 44500
$ws.Range("L116").Value = 44500
$ws.Range("N116").Value = -53678
$ws.Range("H126").Value = 7785.231
$ws.Range("I126").Value = 9030.799999999999
$ws.Range("K126").Value = 27092.4
$ws.Range("M126").Value = -24622.4
